# Daily attendance processing - reorder "Recorded By" entries in column G
# Change "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Column G is the "Recorded By" column
$colIndex = 7

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $colIndex)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
